# AnyTurnsLeft angepasst + BugFixing
#
# Funktion müsste performanter sein und funktioniert prinzipiell (ohne
# ausgiebiges Testen). Ein Bug wurde mit Daniel gefixed, Steine wurden
# übermalt bei "rechts", "runter", "links" Buttons.
#
# -> Adds a new hour-log entry (02.06.2016, 4h) to the Stundenliste and
#    extends the running total accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New log entry in row 26 -------------------------------------------------
$ws.Cells.Item(26, 1).Value = Get-Date -Year 2016 -Month 6 -Day 2 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(26, 2).Value = 4
$ws.Cells.Item(26, 3).Value = "Funktion die überprüft… überarbeitet (Performanter, noch nicht ausgiebig getestet), BugFixing"

# Give the new row the same boxed-in look (thin left/right border) that the
# other filled-in rows in the sheet carry.
foreach ($addr in @("B26", "C26")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}

# --- Extend the hours total to include the new row --------------------------
$ws.Range("B28").Formula = "=SUM(B3:B26)"

# --- Update the on-screen selection -----------------------------------------
$ws.Range("A27").Select()
